$wb = $excel.ActiveWorkbook

# New event row (2024-05-16, the earliest-dated event) that needs to be
# inserted as the new row 2 (right after the header row) on both the
# "展览" (exhibition) sheet and the combined "全部类型" (all types) sheet.
# Inserting the row pushes every following row down by one.
$newB = "2024-05-16"
$newC = "南宁·2024年中国(华南)动漫超级订货会暨动漫实业发展大会"
$newD = "月秋路9号 南宁五象山庄"
$newE = "2024.05.16 09:00-05.16 18:00"
$newF = 1
$newG = 100
$newH = "https://show.bilibili.com/platform/detail.html?id=85180"
$newI = "//i1.hdslb.com/bfs/openplatform/202405/vi1DMzSX1714564480040.jpeg"

# A handful of "想去人数" (interested-attendee count) values ticked up
# slightly between the two snapshots, keyed by the (stable) event name so
# the lookup survives the row-shift caused by the insert above.
$fTweaks = @{
    "南宁·AP动漫游戏嘉年华"     = 2485
    "南宁·布谷鸟动漫展4th"      = 1809
    "南宁·恋与深空only"         = 344
    "南宁·小蜜蜂动漫嘉年华2.0"  = 101
    "南宁·AB动漫游戏嘉年华"     = 861
}

function Update-ExpoSheet($ws) {
    # Push all data rows down by one, opening up row 2 for the new event.
    $ws.Rows.Item(2).Insert()

    # Fill the freshly-inserted row 2. B2 ("2024-05-16") looks like a date,
    # so force the cell to Text first to stop Excel silently reinterpreting
    # it as a date serial, then drop the formatting back to the workbook's
    # default ("Normal") style so the cell matches its unstyled neighbours
    # (B3, B4, ...) once the literal string is safely stored.
    $ws.Cells.Item(2, 2).NumberFormat = "@"
    $ws.Cells.Item(2, 2).Value = $newB
    $ws.Cells.Item(2, 2).Style = "Normal"

    $ws.Cells.Item(2, 3).Value = $newC
    $ws.Cells.Item(2, 4).Value = $newD
    $ws.Cells.Item(2, 5).Value = $newE
    $ws.Cells.Item(2, 6).Value = $newF
    $ws.Cells.Item(2, 7).Value = $newG
    $ws.Cells.Item(2, 8).Value = $newH
    $ws.Cells.Item(2, 9).Value = $newI

    # Column A holds a bold/centered/bordered running index (style matches
    # the rest of the column); reproduce that formatting explicitly so the
    # new cell resolves to the same shared style as its neighbours.
    $ws.Cells.Item(2, 1).Value = 1
    $ws.Cells.Item(2, 1).Font.Bold = $true
    $ws.Cells.Item(2, 1).HorizontalAlignment = -4108
    $ws.Cells.Item(2, 1).VerticalAlignment = -4160
    $ws.Cells.Item(2, 1).Borders.LineStyle = 1

    # Renumber the running index (column A) for every remaining data row.
    $used = $ws.UsedRange
    $lastRow = $used.Row + $used.Rows.Count - 1
    for ($r = 3; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }

    # Apply the small "想去人数" upticks, matched by event name so the
    # lookup is independent of the row shift. (.Value2 is used for the
    # read-back since .Value does not reliably marshal string content here.)
    for ($r = 3; $r -le $lastRow; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($fTweaks.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $fTweaks[$name]
        }
    }
}

Update-ExpoSheet $wb.Worksheets.Item("展览")
Update-ExpoSheet $wb.Worksheets.Item("全部类型")
